# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The existing sheet ends at column AC (A1:AC61); we append three new
# columns - AD "Wins", AE "Losses", AF "Ties" - with a header row (row 1)
# and a constant record (63 wins, 99 losses, 0 ties) repeated for every
# player row (rows 2-61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 61

# --- Header row -------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header cells (bold, centered, bordered)
# by copying A1's formatting onto the new header cells without touching
# the text we just set.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Range("AD$r").Value = 63
    $ws.Range("AE$r").Value = 99
    $ws.Range("AF$r").Value = 0
}
